$wb = $excel.ActiveWorkbook

$wsReservoirs = $wb.Worksheets.Item("Reservoirs")
$wsCoord = $wb.Worksheets.Item("CoordinatedOps")

# --- Add a new row of data (another year) to CoordinatedOps: row 5 ---
# Copy formatting (style/number format) from the row above, then fill in values.
$wsCoord.Cells.Item(4, 1).Copy()
$wsCoord.Cells.Item(5, 1).PasteSpecial(-4122)   # xlPasteFormats

$wsCoord.Cells.Item(5, 1).Value = "1/1/2021"
$wsCoord.Cells.Item(5, 2).Value = "NaN"
$wsCoord.Cells.Item(5, 3).Value = "NaN"
$wsCoord.Cells.Item(5, 4).Value = "NaN"
$wsCoord.Cells.Item(5, 5).Value = "NaN"
$wsCoord.Cells.Item(5, 6).Value = 3684.76
$wsCoord.Cells.Item(5, 7).Value = 1207.9000000000001

# --- Update selections so the Reservoirs sheet selection moves off the
#     full-column pick and onto C25 ---
$wsReservoirs.Range("C25").Select() | Out-Null

# --- CoordinatedOps becomes the active/selected tab, with F14 selected ---
$wsCoord.Activate()
$wsCoord.Range("F14").Select() | Out-Null
